$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update the three multi-run summary rows near the bottom of the table
# (rows 44-46 in the original 46-row table). These collapse a run-per-tab
# value list down to a single summary value.
$t.Cell(44, 1).Range.Text = "99.96"
$t.Cell(45, 1).Range.Text = "0.03"
$t.Cell(46, 1).Range.Text = "63"

# --- Remove the three trailing rows that are folded into row 9's new value
# (originally rows 10, 11, 12: 0.00018 / 0.00019 / 0.00927).
$t.Rows.Item(12).Delete()
$t.Rows.Item(11).Delete()
$t.Rows.Item(10).Delete()

# --- Update rows 6-9 (0.00030 / 0.00014 / 0.00004 / 0.00016) to their new values.
$t.Cell(6, 1).Range.Text = "0.00024"
$t.Cell(7, 1).Range.Text = "0.00027"
$t.Cell(8, 1).Range.Text = "0.00038"
$t.Cell(9, 1).Range.Text = "0.02714"

# --- Row 4 (62 -> 144), then insert three new rows after it.
$t.Cell(4, 1).Range.Text = "144"

$newRow1 = $t.Rows.Add($t.Rows.Item(5))
$newRow1.Cells.Item(1).Range.Text = "0.00003"

$newRow2 = $t.Rows.Add($t.Rows.Item(6))
$newRow2.Cells.Item(1).Range.Text = "0.00059"

$newRow3 = $t.Rows.Add($t.Rows.Item(7))
$newRow3.Cells.Item(1).Range.Text = "0.00017"

# --- Rows 1-3 (99.96 / 0.03 / 63) all become "0M".
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
